$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Förändrad" (changed) date column (C) for rows 2-7
# from serial 45178 (2023-09-09) to serial 45179 (2023-09-10)
$ws.Range("C2:C7").Value = 45179
